$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H129").Value = 929.51666
$ws.Range("J129").Value = 968.0702
$ws.Range("L129").Value = 2904.2106
$ws.Range("N129").Value = -12904.2106
$ws.Range("H132").Value = 877.3521
$ws.Range("I132").Value = 737.9077
$ws.Range("J132").Value = 2388
$ws.Range("K132").Value = 2213.7231
$ws.Range("L132").Value = 7164
$ws.Range("M132").Value = 316.2768999999998
$ws.Range("N132").Value = -12224

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 6998.0586
$ws.Range("I45").Value = 7806.8667
$ws.Range("J45").Value = 932
$ws.Range("K45").Value = 7806.8667
$ws.Range("L45").Value = 932
$ws.Range("M45").Value = -7429.8667
$ws.Range("N45").Value = -1686
$ws.Range("H61").Value = 2625.66
$ws.Range("I61").Value = 1867.7587
$ws.Range("J61").Value = 3672.2856
$ws.Range("K61").Value = 1867.7587
$ws.Range("L61").Value = 3672.2856
$ws.Range("M61").Value = -1655.7587
$ws.Range("N61").Value = -4096.2856
$ws.Range("H74").Value = 1401.3125
$ws.Range("I74").Value = 1220
$ws.Range("J74").Value = 1800.2
$ws.Range("K74").Value = 1220
$ws.Range("L74").Value = 1800.2
$ws.Range("M74").Value = -346
$ws.Range("N74").Value = -3548.2
$ws.Range("H77").Value = 1401.3125
$ws.Range("I77").Value = 1220
$ws.Range("J77").Value = 1800.2
$ws.Range("K77").Value = 6100
$ws.Range("L77").Value = 9001
$ws.Range("M77").Value = -1732
$ws.Range("N77").Value = -17737
$ws.Range("H102").Value = 2316804
$ws.Range("I102").Value = 2647511.8
$ws.Range("J102").Value = 1850
$ws.Range("K102").Value = 2647511.8
$ws.Range("L102").Value = 1850
$ws.Range("M102").Value = -2645889.8
$ws.Range("N102").Value = -5094
$ws.Range("H110").Value = 1841.069
$ws.Range("I110").Value = 2076.2104
$ws.Range("J110").Value = 1394.3
$ws.Range("K110").Value = 2076.2104
$ws.Range("L110").Value = 1394.3
$ws.Range("M110").Value = -31.21039999999994
$ws.Range("N110").Value = -5484.3
$ws.Range("H122").Value = 1426395.9
$ws.Range("I122").Value = 1426395.9
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4279187.699999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4276737.699999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 10006235
$ws.Range("I132").Value = 3083.25
$ws.Range("J132").Value = 16675002
$ws.Range("K132").Value = 9249.75
$ws.Range("L132").Value = 50025006
$ws.Range("M132").Value = -6719.75
$ws.Range("N132").Value = -50030066
$ws.Range("H136").Value = 2625.66
$ws.Range("I136").Value = 1867.7587
$ws.Range("J136").Value = 3672.2856
$ws.Range("K136").Value = 5603.2761
$ws.Range("L136").Value = 11016.8568
$ws.Range("M136").Value = -3053.2761
$ws.Range("N136").Value = -16116.8568

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 931.4666999999999
$ws.Range("I107").Value = 792.1429000000001
$ws.Range("J107").Value = 1256.5555
$ws.Range("K107").Value = 792.1429000000001
$ws.Range("L107").Value = 1256.5555
$ws.Range("M107").Value = 1127.8571
$ws.Range("N107").Value = -5096.5555
$ws.Range("H134").Value = 2573.4814
$ws.Range("I134").Value = 2179.0908
$ws.Range("J134").Value = 2844.625
$ws.Range("K134").Value = 6537.2724
$ws.Range("L134").Value = 8533.875
$ws.Range("M134").Value = -4002.2724
$ws.Range("N134").Value = -13603.875

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2252.45
$ws.Range("I16").Value = 2151.375
$ws.Range("J16").Value = 2319.8333
$ws.Range("K16").Value = 2151.375
$ws.Range("L16").Value = 2319.8333
$ws.Range("M16").Value = -1864.375
$ws.Range("N16").Value = -2893.8333
$ws.Range("H105").Value = 1593.5
$ws.Range("I105").Value = 1703.625
$ws.Range("J105").Value = 1299.8334
$ws.Range("K105").Value = 1703.625
$ws.Range("L105").Value = 1299.8334
$ws.Range("M105").Value = 43.375
$ws.Range("N105").Value = -4793.8334
$ws.Range("H113").Value = 2252.45
$ws.Range("I113").Value = 2151.375
$ws.Range("J113").Value = 2319.8333
$ws.Range("K113").Value = 2151.375
$ws.Range("L113").Value = 2319.8333
$ws.Range("M113").Value = 18.625
$ws.Range("N113").Value = -6659.8333
$ws.Range("H132").Value = 1915.037
$ws.Range("I132").Value = 1533.4706
$ws.Range("J132").Value = 2563.7
$ws.Range("K132").Value = 4600.4118
$ws.Range("L132").Value = 7691.099999999999
$ws.Range("M132").Value = -2070.4118
$ws.Range("N132").Value = -12751.1
$ws.Range("H134").Value = 2374.5908
$ws.Range("I134").Value = 2438.0264
$ws.Range("J134").Value = 1972.8334
$ws.Range("K134").Value = 7314.0792
$ws.Range("L134").Value = 5918.5002
$ws.Range("M134").Value = -4779.0792
$ws.Range("N134").Value = -10988.5002

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 570.8182
$ws.Range("J107").Value = 741.0714
$ws.Range("L107").Value = 2223.2142
$ws.Range("N107").Value = -6063.2142
$ws.Range("H122").Value = 2489.1372
$ws.Range("I122").Value = 485.4
$ws.Range("J122").Value = 2706.9348
$ws.Range("K122").Value = 4368.599999999999
$ws.Range("L122").Value = 24362.4132
$ws.Range("M122").Value = -1918.599999999999
$ws.Range("N122").Value = -29262.4132
$ws.Range("H130").Value = 4185.484
$ws.Range("I130").Value = 1815
$ws.Range("J130").Value = 4348.9653
$ws.Range("K130").Value = 5445
$ws.Range("L130").Value = 13046.8959
$ws.Range("M130").Value = -425
$ws.Range("N130").Value = -23086.8959

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 40001270
$ws.Range("I113").Value = 52632456
$ws.Range("J113").Value = 2528.8333
$ws.Range("K113").Value = 52632456
$ws.Range("L113").Value = 2528.8333
$ws.Range("M113").Value = -52630286
$ws.Range("N113").Value = -6868.8333
$ws.Range("H122").Value = 39440196
$ws.Range("I122").Value = 56045056
$ws.Range("K122").Value = 168135168
$ws.Range("M122").Value = -168132718
$ws.Range("H123").Value = 24523.238
$ws.Range("J123").Value = 25299.4
$ws.Range("L123").Value = 25299.4
$ws.Range("N123").Value = -30199.4
$ws.Range("H126").Value = 5440.2
$ws.Range("I126").Value = 7411.1763
$ws.Range("J126").Value = 2862.7693
$ws.Range("K126").Value = 22233.5289
$ws.Range("L126").Value = 8588.3079
$ws.Range("M126").Value = -19763.5289
$ws.Range("N126").Value = -13528.3079
$ws.Range("H132").Value = 3550.8076
$ws.Range("I132").Value = 5330.6665
$ws.Range("J132").Value = 2608.5293
$ws.Range("K132").Value = 15991.9995
$ws.Range("L132").Value = 7825.5879
$ws.Range("M132").Value = -13461.9995
$ws.Range("N132").Value = -12885.5879

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 265.45
$ws.Range("I55").Value = 250.66667
$ws.Range("J55").Value = 277.54544
$ws.Range("K55").Value = 250.66667
$ws.Range("L55").Value = 277.54544
$ws.Range("M55").Value = -77.66667000000001
$ws.Range("N55").Value = -623.54544
$ws.Range("H61").Value = 1334.85
$ws.Range("I61").Value = 1318.5625
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 1318.5625
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -1116.5625
$ws.Range("N61").Value = -1804
$ws.Range("H113").Value = 1334.85
$ws.Range("I113").Value = 1318.5625
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1318.5625
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 851.4375
$ws.Range("N113").Value = -5740
$ws.Range("H122").Value = 4291762.5
$ws.Range("I122").Value = 4470843
$ws.Range("J122").Value = 3336666.8
$ws.Range("K122").Value = 13412529
$ws.Range("L122").Value = 10010000.4
$ws.Range("M122").Value = -13410079
$ws.Range("N122").Value = -10014900.4
$ws.Range("H132").Value = 20837596
$ws.Range("I132").Value = 41669860
$ws.Range("J132").Value = 5331.625
$ws.Range("K132").Value = 125009580
$ws.Range("L132").Value = 15994.875
$ws.Range("M132").Value = -125007050
$ws.Range("N132").Value = -21054.875
$ws.Range("H136").Value = 3639.9033
$ws.Range("I136").Value = 1870.6731
$ws.Range("J136").Value = 12839.9
$ws.Range("K136").Value = 5612.0193
$ws.Range("L136").Value = 38519.7
$ws.Range("M136").Value = -3062.0193
$ws.Range("N136").Value = -43619.7

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 1371.2667
$ws.Range("I96").Value = 1390.6428
$ws.Range("K96").Value = 1390.6428
$ws.Range("M96").Value = -17.64280000000008
$ws.Range("H107").Value = 166667660
$ws.Range("I107").Value = 200000860
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 600002580
$ws.Range("L107").Value = 4800
$ws.Range("M107").Value = -600000660
$ws.Range("N107").Value = -8640
$ws.Range("H113").Value = 871.85187
$ws.Range("I113").Value = 774.61536
$ws.Range("K113").Value = 2323.84608
$ws.Range("M113").Value = -153.8460800000003
$ws.Range("H122").Value = 1260.2
$ws.Range("I122").Value = 1000.5
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 3001.5
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -551.5
$ws.Range("N122").Value = -9200.0002
$ws.Range("H132").Value = 1531.75
$ws.Range("I132").Value = 1082.7931
$ws.Range("J132").Value = 2399.7334
$ws.Range("K132").Value = 3248.379300000001
$ws.Range("L132").Value = 7199.2002
$ws.Range("M132").Value = -718.3793000000005
$ws.Range("N132").Value = -12259.2002
$ws.Range("H136").Value = 2656.081
$ws.Range("I136").Value = 3411.5625
$ws.Range("J136").Value = 2080.476
$ws.Range("K136").Value = 10234.6875
$ws.Range("L136").Value = 6241.428
$ws.Range("M136").Value = -7684.6875
$ws.Range("N136").Value = -11341.428
